# Applies the recorded edit:
#  1. Re-colours the deck's live theme (ppt/theme/theme2.xml, used by the
#     SlideMaster / all slides) from the "Integral" (Red Violet) palette to
#     the standard "Office Theme" palette.
#  2. Re-applies the default table style (tableStyleId) on the three tables
#     that had gone through the table-style gallery.

$p = $ppt.ActivePresentation

function RGBFromHex($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette: built-in "Office Theme" colour scheme.
$officeColors = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = RGBFromHex($officeColors[$i - 1])
}

# Re-apply the (now differently-GUID'd) default table style to every table
# in the deck that is still wearing the old style id.
$oldStyleId = "{EAC139C7-D7D8-4AC5-9ED2-49A51604CE8E}"
$newStyleId = "{F15D04FD-3D67-4645-BC75-526018EEF15D}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
